$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 107 (hunk 0)
$ws.Range("H107").Value = 1225
$ws.Range("I107").Value = 1225
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1225
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 695
$ws.Range("N107").Value = ""

# row 129 (hunk 1)
$ws.Range("H129").Value = 773.8182
$ws.Range("I129").Value = 559.1818
$ws.Range("J129").Value = 988.4545000000001
$ws.Range("K129").Value = 1677.5454
$ws.Range("L129").Value = 2965.3635
$ws.Range("M129").Value = 3322.4546
$ws.Range("N129").Value = -12965.3635

# row 132 (hunk 2)
$ws.Range("H132").Value = 723292.3
$ws.Range("I132").Value = 2532.85
$ws.Range("J132").Value = 6128988
$ws.Range("K132").Value = 7598.549999999999
$ws.Range("L132").Value = 18386964
$ws.Range("M132").Value = -5068.549999999999

# row 135 (hunk 3)
$ws.Range("H135").Value = 28937.973
$ws.Range("I135").Value = 37334
$ws.Range("J135").Value = 2817
$ws.Range("K135").Value = 336006
$ws.Range("L135").Value = 25353
$ws.Range("M135").Value = -333471
$ws.Range("N135").Value = -30423

# row 137 (hunk 4)
$ws.Range("H137").Value = 2130439.2
$ws.Range("I137").Value = 4001961.5
$ws.Range("J137").Value = 3709.2273
$ws.Range("K137").Value = 12005884.5
$ws.Range("L137").Value = 11127.6819
$ws.Range("M137").Value = -12003334.5


# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 32 (hunk 5)
$ws.Range("H32").Value = 4559.13
$ws.Range("I32").Value = 4184.811
$ws.Range("J32").Value = 7928
$ws.Range("K32").Value = 4184.811
$ws.Range("L32").Value = 7928
$ws.Range("M32").Value = -3897.811
$ws.Range("N32").Value = -8502

# row 82 (hunk 6)
$ws.Range("H82").Value = 40181
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 40181
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 40181
$ws.Range("N82").Value = -40903

# row 85 (hunk 7)
$ws.Range("H85").Value = 40181
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 40181
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 40181
$ws.Range("N85").Value = -42677

# row 132 (hunk 8)
$ws.Range("H132").Value = 13214820
$ws.Range("I132").Value = 15657956
$ws.Range("J132").Value = 184766.33
$ws.Range("K132").Value = 46973868
$ws.Range("L132").Value = 554298.99
$ws.Range("M132").Value = -46971338
$ws.Range("N132").Value = -559358.99


# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31 (hunk 9)
$ws.Range("H31").Value = 346121.38
$ws.Range("I31").Value = 79342.92
$ws.Range("J31").Value = 503763.2
$ws.Range("K31").Value = 79342.92
$ws.Range("L31").Value = 503763.2
$ws.Range("M31").Value = -79047.92
$ws.Range("N31").Value = -504353.2

# row 34 (hunk 10)
$ws.Range("H34").Value = 346121.38
$ws.Range("I34").Value = 79342.92
$ws.Range("J34").Value = 503763.2
$ws.Range("K34").Value = 79342.92
$ws.Range("L34").Value = 503763.2
$ws.Range("M34").Value = -79140.92
$ws.Range("N34").Value = -504167.2

# row 134 (hunk 11)
$ws.Range("H134").Value = 40910.258
$ws.Range("I134").Value = 617.0769
$ws.Range("J134").Value = 78325.36
$ws.Range("K134").Value = 1851.2307
$ws.Range("L134").Value = 234976.08
$ws.Range("M134").Value = 683.7692999999999


# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 6 (hunk 12)
$ws.Range("H6").Value = 306.7143
$ws.Range("I6").Value = 57.833332
$ws.Range("J6").Value = 1800
$ws.Range("K6").Value = 173.499996
$ws.Range("L6").Value = 5400
$ws.Range("M6").Value = -60.49999600000001
$ws.Range("N6").Value = -5626

# row 18 (hunk 13)
$ws.Range("H18").Value = 186.125
$ws.Range("I18").Value = 136.78572
$ws.Range("J18").Value = 531.5
$ws.Range("K18").Value = 410.35716
$ws.Range("L18").Value = 1594.5
$ws.Range("M18").Value = -241.35716
$ws.Range("N18").Value = -1932.5

# row 25 (hunk 14)
$ws.Range("H25").Value = 500
$ws.Range("I25").Value = 500
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 1500
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -1331

# row 30 (hunk 15)
$ws.Range("H30").Value = 500
$ws.Range("I30").Value = 500
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 1500
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -1398

# row 54 (hunk 16)
$ws.Range("H54").Value = 3000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 3000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 9000
$ws.Range("N54").Value = -10118

# row 56 (hunk 17)
$ws.Range("H56").Value = 179833.31
$ws.Range("I56").Value = 179833.31
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 179833.31
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -179303.31

# row 70 (hunk 18)
$ws.Range("H70").Value = 2355.7334
$ws.Range("I70").Value = 1048
$ws.Range("J70").Value = 3500
$ws.Range("K70").Value = 3144
$ws.Range("L70").Value = 10500
$ws.Range("M70").Value = -2829

# row 73 (hunk 19)
$ws.Range("H73").Value = 2355.7334
$ws.Range("I73").Value = 1048
$ws.Range("J73").Value = 3500
$ws.Range("K73").Value = 3144
$ws.Range("L73").Value = 10500
$ws.Range("M73").Value = -2052

# row 75 (hunk 20)
$ws.Range("H75").Value = 4547
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 4547
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 13641
$ws.Range("N75").Value = -15637

# row 76 (hunk 21)
$ws.Range("H76").Value = 3091.25
$ws.Range("I76").Value = 915
$ws.Range("J76").Value = 3816.6667
$ws.Range("K76").Value = 2745
$ws.Range("L76").Value = 11450.0001
$ws.Range("M76").Value = -2362
$ws.Range("N76").Value = -12216.0001

# row 78 (hunk 22)
$ws.Range("H78").Value = 4547
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 4547
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 40923
$ws.Range("N78").Value = -50907

# row 79 (hunk 23)
$ws.Range("H79").Value = 3091.25
$ws.Range("I79").Value = 915
$ws.Range("J79").Value = 3816.6667
$ws.Range("K79").Value = 2745
$ws.Range("L79").Value = 11450.0001
$ws.Range("M79").Value = -1419
$ws.Range("N79").Value = -14102.0001

# row 88 (hunk 24)
$ws.Range("H88").Value = 3000
$ws.Range("I88").Value = 3000
$ws.Range("J88").Value = 3000
$ws.Range("K88").Value = 9000
$ws.Range("L88").Value = 9000
$ws.Range("M88").Value = -8572
$ws.Range("N88").Value = -9856

# row 91 (hunk 25)
$ws.Range("H91").Value = 3000
$ws.Range("I91").Value = 3000
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 9000
$ws.Range("L91").Value = 9000
$ws.Range("M91").Value = -7518
$ws.Range("N91").Value = -11964

# row 108 (hunk 26)
$ws.Range("H108").Value = 431.5
$ws.Range("I108").Value = 431.5
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 1294.5
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = 1585.5

# row 114 (hunk 27)
$ws.Range("H114").Value = 10417192
$ws.Range("I114").Value = 458.91666
$ws.Range("J114").Value = 16667232
$ws.Range("K114").Value = 1376.74998
$ws.Range("L114").Value = 50001696
$ws.Range("M114").Value = 1877.25002
$ws.Range("N114").Value = -50008204

# row 122 (hunk 28)
$ws.Range("H122").Value = 1025.9667
$ws.Range("I122").Value = 269.25
$ws.Range("J122").Value = 1301.1364
$ws.Range("K122").Value = 2423.25
$ws.Range("L122").Value = 11710.2276
$ws.Range("M122").Value = 26.75
$ws.Range("N122").Value = -16610.2276

# row 131 (hunk 29)
$ws.Range("H131").Value = 992.75
$ws.Range("I131").Value = 500
$ws.Range("J131").Value = 1063.1428
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 3189.4284
$ws.Range("M131").Value = 3540
$ws.Range("N131").Value = -13269.4284

# row 132 (hunk 30)
$ws.Range("H132").Value = 3525.5
$ws.Range("I132").Value = 2481.6
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 22334.4
$ws.Range("L132").Value = 36000
$ws.Range("M132").Value = -19804.4
$ws.Range("N132").Value = -41060


# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 123 (hunk 31)
$ws.Range("H123").Value = 23087
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 23087
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 23087
$ws.Range("N123").Value = -27987


# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 7 (hunk 32)
$ws.Range("H7").Value = 125002000
$ws.Range("I7").Value = 125002000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 125002000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -125001888
$ws.Range("N7").Value = ""

# row 40 (hunk 33)
$ws.Range("H40").Value = 2454
$ws.Range("I40").Value = 2108
$ws.Range("J40").Value = 2800
$ws.Range("K40").Value = 2108
$ws.Range("L40").Value = 2800
$ws.Range("M40").Value = -1972
$ws.Range("N40").Value = -3072

# row 46 (hunk 34)
$ws.Range("H46").Value = 2525942.8
$ws.Range("I46").Value = 4329430
$ws.Range("J46").Value = 1060.4
$ws.Range("K46").Value = 4329430
$ws.Range("L46").Value = 1060.4
$ws.Range("M46").Value = -4329242
$ws.Range("N46").Value = -1436.4

# row 126 (hunk 35)
$ws.Range("H126").Value = 125002000
$ws.Range("I126").Value = 125002000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 375006000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -375003530
$ws.Range("N126").Value = ""

# row 136 (hunk 36)
$ws.Range("H136").Value = 38856.215
$ws.Range("I136").Value = 23546.184
$ws.Range("J136").Value = 146026.42
$ws.Range("K136").Value = 70638.552
$ws.Range("L136").Value = 438079.26
$ws.Range("M136").Value = -68088.552
$ws.Range("N136").Value = -443179.26


# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 122 (hunk 37)
$ws.Range("H122").Value = 1082.091
$ws.Range("I122").Value = 955
$ws.Range("J122").Value = 1304.5
$ws.Range("K122").Value = 2865
$ws.Range("L122").Value = 3913.5
$ws.Range("M122").Value = -415
$ws.Range("N122").Value = -8813.5

